$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "general" (A/B columns) and "specific" (I/J columns) hybrid-performance
# blocks in rows 16-20 already carried a classifier label in column A/I
# (J48, KNN, NB, RF, SMO). Label each of those rows with the corresponding
# "Hybrid General" / "Hybrid specific" tag in column B/J, matching the
# pattern already used for rows 9-13 (column B/J = "DDosType").
foreach ($row in 16..20) {
    $ws.Cells.Item($row, 2).Value = "Hybrid General"   # column B
    $ws.Cells.Item($row, 10).Value = "Hybrid specific" # column J
}

# Move the active selection to I16 (was L17).
$ws.Range("I16").Select()
